$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "29.495.35"
$ws.Cells.Item(2,5).Value = "  +0.85%  "

# Row 3
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = "1.975.10"
$ws.Cells.Item(3,5).Value = "  +3.97%  "

# Row 4
$ws.Cells.Item(4,5).Value = "  +0.19%  "

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "327.09"
$ws.Cells.Item(5,5).Value = "  +0.24%  "

# Row 6
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "1.004"
$ws.Cells.Item(6,5).Value = "  +0.12%  "

# Row 7
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "0.4666"
$ws.Cells.Item(7,5).Value = "  +0.80%  "

# Row 8
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "0.3922"
$ws.Cells.Item(8,5).Value = "  -0.08%  "

# Row 9
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "46.25"
$ws.Cells.Item(9,5).Value = "  -1.08%  "

# Row 10
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "0.07955"
$ws.Cells.Item(10,5).Value = "  +0.97%  "

# Row 11
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.9909"
$ws.Cells.Item(11,5).Value = "  +0.19%  "

# Row 12
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "22.77"
$ws.Cells.Item(12,5).Value = "  +4.11%  "

# Row 13
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "1.993.94"
$ws.Cells.Item(13,5).Value = "  +4.60%  "

# Row 14
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "7.188"
$ws.Cells.Item(14,5).Value = "  +1.69%  "

# Row 15
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "5.844"
$ws.Cells.Item(15,5).Value = "  +1.84%  "

# Row 16
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "0.07075"
$ws.Cells.Item(16,5).Value = "  +1.20%  "

# Row 17
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "87.79"
$ws.Cells.Item(17,5).Value = "  -0.70%  "

# Row 18
$ws.Cells.Item(18,5).Value = "  +0.07%  "

# Row 19
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "0.000009953"
$ws.Cells.Item(19,5).Value = "  -0.15%  "

# Row 20
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "17.30"
$ws.Cells.Item(20,5).Value = "  +1.39%  "

# Row 21
$ws.Cells.Item(21,5).Value = "  +0.19%  "

# Row 22
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "29.492.13"
$ws.Cells.Item(22,5).Value = "  +0.79%  "

# Row 23
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "0.5070"
$ws.Cells.Item(23,5).Value = "  +6.09%  "

# Row 24
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "5.551"
$ws.Cells.Item(24,5).Value = "  +4.54%  "

# Row 25
$ws.Cells.Item(25,5).Value = "  +0.51%  "

# Row 26
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "2.209.46"
$ws.Cells.Item(26,5).Value = "  +3.19%  "

# Row 27
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "2.110"
$ws.Cells.Item(27,5).Value = "  +0.42%  "

# Row 28
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "158.45"
$ws.Cells.Item(28,5).Value = "  +1.67%  "

# Row 29
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "19.52"
$ws.Cells.Item(29,5).Value = "  +0.68%  "

# Row 30
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "5.774"
$ws.Cells.Item(30,5).Value = "  -4.16%  "

# Row 31
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "119.59"
$ws.Cells.Item(31,5).Value = "  +0.83%  "

# Row 32
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = "1.909"
$ws.Cells.Item(32,5).Value = "  +0.71%  "

# Row 33
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "0.09406"
$ws.Cells.Item(33,5).Value = "  +0.58%  "

# Row 34
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "0.8928"
$ws.Cells.Item(34,5).Value = "  -1.45%  "

# Row 35
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = "5.239"
$ws.Cells.Item(35,5).Value = "  -0.07%  "

# Row 36
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "1.323"
$ws.Cells.Item(36,5).Value = "  -0.06%  "

# Row 37
$ws.Cells.Item(37,5).Value = "  -1.55%  "

# Row 38
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "0.05825"
$ws.Cells.Item(38,5).Value = "  +0.80%  "

# Row 39
$ws.Cells.Item(39,5).Value = "  -1.03%  "

# Row 40
$ws.Cells.Item(40,5).Value = "  +0.73%  "

# Row 41
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "7.773"
$ws.Cells.Item(41,5).Value = "  +0.62%  "

# Row 42
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "0.5727"
$ws.Cells.Item(42,5).Value = "  +0.47%  "

# Row 43
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "0.000003089"
$ws.Cells.Item(43,5).Value = "  +48.83%  "

# Row 44
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "0.1798"
$ws.Cells.Item(44,5).Value = "  +0.75%  "

# Row 45
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "9.677"
$ws.Cells.Item(45,5).Value = "  -0.33%  "

# Row 46
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "2.758"
$ws.Cells.Item(46,5).Value = "  +7.06%  "

# Row 47
$ws.Cells.Item(47,2).Value = "EnergySwap"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "11.81"
$ws.Cells.Item(47,5).Value = "  -1.18%  "

# Row 48
$ws.Cells.Item(48,2).Value = "RenderToken"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "2.206"
$ws.Cells.Item(48,5).Value = "  +1.74%  "

# Row 49
$ws.Cells.Item(49,2).Value = "Decentraland"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "0.5353"
$ws.Cells.Item(49,5).Value = "  +0.01%  "

# Row 50
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "0.06929"
$ws.Cells.Item(50,5).Value = "  -1.47%  "

# Row 51
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "1.831"
$ws.Cells.Item(51,5).Value = "  -1.02%  "
